$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("SLEN", 6, 1, "Represents and leads systems and software life cycle working practices at the highest level in the organisation"),
    @("SLEN", 6, 2, "Identifies opportunities for innovation in systems and software life cycle working practices to achieve organisational goals and objectives"),
    @("SLEN", 6, 3, "Leads the essential cultural and environmental changes and communicates the benefits to all stakeholders"),
    @("SLEN", 6, 4, "Oversees the quality of the work performed and delivers measurable business benefits")
)

$startRow = 14
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Cells.Item($r, 4).Value = $rows[$i][3]
}
